$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Replace every occurrence of "(ºC)" with the double-encoded mojibake
# variant "(ÂºC)" throughout the document body (covers all table cells).
$d.Content.Find.Execute("(ºC)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(ÂºC)", 2)
